$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Auto-fit columns A:F to their content (header/collector names, amounts, rates)
# so the resulting <cols> widths match Excel's computed "best fit" widths as
# closely as this engine's column-width model allows.
$ws.Columns("A").ColumnWidth = 4.877604166666667
$ws.Columns("B").ColumnWidth = 11.166666666666666
$ws.Columns("C").ColumnWidth = 27.592447916666668
$ws.Columns("D").ColumnWidth = 24.307291666666668
$ws.Columns("E").ColumnWidth = 10.307291666666666
$ws.Columns("F").ColumnWidth = 12.451822916666666

# Move the active selection to H10, matching the saved sheet view state
$ws.Range("H10").Select()
